{"js": "const replacements = [\n  [\"2025-01-16 Thursday\", \"2025-01-17 Friday\"],\n  [\"389\u00f78=\", \"933\u00f79=\"],\n  [\"730\u00f75=\", \"443\u00f79=\"],\n  [\"160\u00f74=\", \"700\u00f76=\"],\n  [\"727\u00f79=\", \"579\u00f79=\"],\n  [\"772\u00f74=\", \"866\u00f75=\"],\n  [\"118\u00f75=\", \"825\u00f76=\"],\n  [\"609\u00f72=\", \"309\u00f76=\"],\n  [\"406\u00f72=\", \"410\u00f78=\"],\n  [\"294\u00f78=\", \"154\u00f75=\"],\n  [\"414\u00f76=\", \"472\u00f78=\"],\n  [\"650\u00f79=\", \"763\u00f74=\"],\n  [\"214\u00f73=\", \"738\u00f72=\"],\n  [\"654\u00f79=\", \"355\u00f75=\"],\n  [\"507\u00f72=\", \"321\u00f78=\"],\n  [\"441\u00f78=\", \"970\u00f78=\"],\n  [\"574\u00f73=\", \"478\u00f73=\"],\n  [\"538\u00f76=\", \"956\u00f77=\"],\n  [\"213\u00f76=\", \"978\u00f78=\"],\n  [\"418\u00f78=\", \"788\u00f73=\"],\n  [\"533\u00f76=\", \"830\u00f74=\"],\n  [\"390\u00f73=\", \"782\u00f78=\"],\n  [\"625\u00f72=\", \"643\u00f79=\"],\n  [\"705\u00f74=\", \"136\u00f79=\"],\n  [\"859\u00f76=\", \"451\u00f77=\"],\n  [\"186\u00f78=\", \"878\u00f79=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"2025-01-16 Thursday\"; New = \"2025-01-17 Friday\" },\n    @{ Old = \"389\u00f78=\"; New = \"933\u00f79=\" },\n    @{ Old = \"730\u00f75=\"; New = \"443\u00f79=\" },\n    @{ Old = \"160\u00f74=\"; New = \"700\u00f76=\" },\n    @{ Old = \"727\u00f79=\"; New = \"579\u00f79=\" },\n    @{ Old = \"772\u00f74=\"; New = \"866\u00f75=\" },\n    @{ Old = \"118\u00f75=\"; New = \"825\u00f76=\" },\n    @{ Old = \"609\u00f72=\"; New = \"309\u00f76=\" },\n    @{ Old = \"406\u00f72=\"; New = \"410\u00f78=\" },\n    @{ Old = \"294\u00f78=\"; New = \"154\u00f75=\" },\n    @{ Old = \"414\u00f76=\"; New = \"472\u00f78=\" },\n    @{ Old = \"650\u00f79=\"; New = \"763\u00f74=\" },\n    @{ Old = \"214\u00f73=\"; New = \"738\u00f72=\" },\n    @{ Old = \"654\u00f79=\"; New = \"355\u00f75=\" },\n    @{ Old = \"507\u00f72=\"; New = \"321\u00f78=\" },\n    @{ Old = \"441\u00f78=\"; New = \"970\u00f78=\" },\n    @{ Old = \"574\u00f73=\"; New = \"478\u00f73=\" },\n    @{ Old = \"538\u00f76=\"; New = \"956\u00f77=\" },\n    @{ Old = \"213\u00f76=\"; New = \"978\u00f78=\" },\n    @{ Old = \"418\u00f78=\"; New = \"788\u00f73=\" },\n    @{ Old = \"533\u00f76=\"; New = \"830\u00f74=\" },\n    @{ Old = \"390\u00f73=\"; New = \"782\u00f78=\" },\n    @{ Old = \"625\u00f72=\"; New = \"643\u00f79=\" },\n    @{ Old = \"705\u00f74=\"; New = \"136\u00f79=\" },\n    @{ Old = \"859\u00f76=\"; New = \"451\u00f77=\" },\n    @{ Old = \"186\u00f78=\"; New = \"878\u00f79=\" }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.Old\n    $find.Replacement.Text = $r.New\n    $find.Execute($r.Old, $false, $false, $false, $false, $false, $true, 1, $false, $r.New, 2) | Out-Null\n}\n"}
